$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = -1.19942097753403
$arr[0,1] = -2.158355674041143
$arr[0,2] = 0.6753625508398458
$arr[0,3] = -0.9937421128663182
$arr[0,4] = 0.02643032487692459
$arr[0,5] = -1.687359132022387
$arr[0,6] = 1.140996241576585
$arr[0,7] = -1.094336358289534
$arr[0,8] = 0.2569986996281282
$arr[0,9] = -0.3644392301887736
$ws.Range("B2:K2").Value2 = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = -2.266970206257284
$arr[0,1] = 0.5667480186237051
$arr[0,2] = -1.102356645082459
$arr[0,3] = -0.08218420733921622
$arr[0,4] = -1.795973664238528
$arr[0,5] = 1.032381709360444
$arr[0,6] = -1.202950890505675
$arr[0,7] = 0.1483841674119874
$arr[0,8] = -0.4730537624049144
$arr[0,9] = 0.02750693478591659
$ws.Range("B3:K3").Value2 = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0.774849739591444
$arr[0,1] = -0.89425492411472
$arr[0,2] = 0.1259175136285228
$arr[0,3] = -1.587871943270789
$arr[0,4] = 1.240483430328183
$arr[0,5] = -0.994849169537936
$arr[0,6] = 0.3564858883797264
$arr[0,7] = -0.2649520414371754
$arr[0,8] = 0.2356086557536556
$arr[0,9] = -0.4322994165924858
$ws.Range("B4:K4").Value2 = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = -0.7777567537409195
$arr[0,1] = 0.2424156840023232
$arr[0,2] = -1.471373772896988
$arr[0,3] = 1.356981600701984
$arr[0,4] = -0.8783509991641355
$arr[0,5] = 0.4729840587535268
$arr[0,6] = -0.148453871063375
$arr[0,7] = 0.352106826127456
$arr[0,8] = -0.3158012462186854
$arr[0,9] = 0.1131997290193177
$ws.Range("B5:K5").Value2 = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 0.2879090979994584
$arr[0,1] = -1.425880358899853
$arr[0,2] = 1.402475014699119
$arr[0,3] = -0.8328575851670005
$arr[0,4] = 0.5184774727506619
$arr[0,5] = -0.1029604570662399
$arr[0,6] = 0.3976002401245912
$arr[0,7] = -0.2703078322215502
$arr[0,8] = 0.1586931430164528
$arr[0,9] = 0.2163646915946629
$ws.Range("B6:K6").Value2 = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = -1.443434480259818
$arr[0,1] = 1.384920893339154
$arr[0,2] = -0.8504117065269649
$arr[0,3] = 0.5009233513906975
$arr[0,4] = -0.1205145784262043
$arr[0,5] = 0.3800461187646267
$arr[0,6] = -0.2878619535815147
$arr[0,7] = 0.1411390216564884
$arr[0,8] = 0.1988105702346985
$arr[0,9] = 0.322776941072984
$ws.Range("B7:K7").Value2 = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 1.473028212290161
$arr[0,1] = -0.7623043875759586
$arr[0,2] = 0.5890306703417038
$arr[0,3] = -0.0324072594751981
$arr[0,4] = 0.4681534377156329
$arr[0,5] = -0.1997546346305085
$arr[0,6] = 0.2292463406074946
$arr[0,7] = 0.2869178891857047
$arr[0,8] = 0.4108842600239903
$arr[0,9] = -0.4825338632108016
$ws.Range("B8:K8").Value2 = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = -0.5837297540881751
$arr[0,1] = 0.7676053038294873
$arr[0,2] = 0.1461673740125855
$arr[0,3] = 0.6467280712034165
$arr[0,4] = -0.02118000114272489
$arr[0,5] = 0.4078209740952782
$arr[0,6] = 0.4654925226734883
$arr[0,7] = 0.5894588935117738
$arr[0,8] = -0.303959229723018
$arr[0,9] = 0.4661714972207444
$ws.Range("B9:K9").Value2 = $arr

$arr = New-Object 'object[,]' 1,10
$arr[0,0] = 1.679632531582137
$arr[0,1] = 1.058194601765235
$arr[0,2] = 1.558755298956066
$arr[0,3] = 0.8908472266099251
$arr[0,4] = 1.319848201847928
$arr[0,5] = 1.377519750426138
$arr[0,6] = 1.501486121264424
$arr[0,7] = 0.608067998029632
$arr[0,8] = 1.378198724973394
$arr[0,9] = 1.11229800409388
$ws.Range("B10:K10").Value2 = $arr

$arr = New-Object 'object[,]' 1,9
$arr[0,0] = 0.1142203657994787
$arr[0,1] = 0.6147810629903097
$arr[0,2] = -0.0531270093558317
$arr[0,3] = 0.3758739658821714
$arr[0,4] = 0.4335455144603815
$arr[0,5] = 0.557511885298667
$arr[0,6] = -0.3359062379361248
$arr[0,7] = 0.4342244890076376
$arr[0,8] = 0.1683237681281231
$ws.Range("B11:J11").Value2 = $arr
$ws.Range("K11").ClearContents()

$arr = New-Object 'object[,]' 1,8
$arr[0,0] = 0.6187489605034189
$arr[0,1] = -0.04915911184272259
$arr[0,2] = 0.3798418633952805
$arr[0,3] = 0.4375134119734906
$arr[0,4] = 0.5614797828117761
$arr[0,5] = -0.3319383404230157
$arr[0,6] = 0.4381923865207467
$arr[0,7] = 0.1722916656412322
$ws.Range("B12:I12").Value2 = $arr
$ws.Range("J12").ClearContents()

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0.08648097832751878
$arr[0,1] = 0.5154819535655218
$arr[0,2] = 0.573153502143732
$arr[0,3] = 0.6971198729820175
$arr[0,4] = -0.1962982502527744
$arr[0,5] = 0.5738324766909881
$arr[0,6] = 0.3079317558114735
$ws.Range("B13:H13").Value2 = $arr
$ws.Range("I13").ClearContents()

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 0.2746757717098572
$arr[0,1] = 0.3323473202880673
$arr[0,2] = 0.4563136911263528
$arr[0,3] = -0.4371044321084391
$arr[0,4] = 0.3330262948353234
$arr[0,5] = 0.06712557395580883
$ws.Range("B14:G14").Value2 = $arr
$ws.Range("H14").ClearContents()

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 0.2870161050359709
$arr[0,1] = 0.4109824758742565
$arr[0,2] = -0.4824356473605354
$arr[0,3] = 0.287695079583227
$arr[0,4] = 0.02179435870371246
$ws.Range("B15:F15").Value2 = $arr
$ws.Range("G15").ClearContents()

$arr = New-Object 'object[,]' 1,4
$arr[0,0] = 0.3441210539382026
$arr[0,1] = -0.5492970692965893
$arr[0,2] = 0.2208336576471732
$arr[0,3] = -0.04506706323234141
$ws.Range("B16:E16").Value2 = $arr
$ws.Range("F16").ClearContents()

$arr = New-Object 'object[,]' 1,3
$arr[0,0] = -0.5788832716533059
$arr[0,1] = 0.1912474552904566
$arr[0,2] = -0.07465326558905801
$ws.Range("B17:D17").Value2 = $arr
$ws.Range("E17").ClearContents()

$arr = New-Object 'object[,]' 1,2
$arr[0,0] = 0.1730967985608157
$arr[0,1] = -0.0928039223186989
$ws.Range("B18:C18").Value2 = $arr
$ws.Range("D18").ClearContents()

$ws.Range("B19").Value2 = -0.1108357465673982
$ws.Range("C19").ClearContents()

$ws.Range("B20").ClearContents()
